# Update the dSF (column F) values per the repulled data / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -1
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = 4
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = -1
$ws.Range("F25").Value = 0
$ws.Range("F28").Value = -2
$ws.Range("F37").Value = 3
$ws.Range("F38").Value = 2
$ws.Range("F39").Value = -4
$ws.Range("F43").Value = -2
$ws.Range("F46").Value = 9
$ws.Range("F47").Value = 2
$ws.Range("F48").Value = 0
